$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-10-25 Friday" "2024-10-26 Saturday"
Replace-Text "735×9=6615" "985×9=8865"
Replace-Text "987×4=3948" "997×5=4985"
Replace-Text "435×5=2175" "297×7=2079"
Replace-Text "156×4=624" "233×9=2097"
Replace-Text "408×8=3264" "917×3=2751"
Replace-Text "289×8=2312" "834×7=5838"
Replace-Text "361×8=2888" "496×9=4464"
Replace-Text "936×4=3744" "695×9=6255"
Replace-Text "750×3=2250" "752×3=2256"
Replace-Text "465×2=930" "287×3=861"
Replace-Text "342×4=1368" "418×8=3344"
Replace-Text "385×4=1540" "701×8=5608"
Replace-Text "840×6=5040" "355×3=1065"
Replace-Text "460×2=920" "237×9=2133"
Replace-Text "204×2=408" "125×2=250"
Replace-Text "279×4=1116" "854×4=3416"
Replace-Text "121×9=1089" "473×2=946"
Replace-Text "538×3=1614" "583×3=1749"
Replace-Text "422×7=2954" "574×8=4592"
Replace-Text "754×6=4524" "610×6=3660"
Replace-Text "160×8=1280" "653×7=4571"
Replace-Text "802×7=5614" "595×4=2380"
Replace-Text "652×9=5868" "288×8=2304"
Replace-Text "899×2=1798" "970×2=1940"
Replace-Text "675×9=6075" "339×2=678"
